# Generate Report for Handoff
# Refreshes the "Latest Handoff Date"/"Latest Handoff Datetime" timestamps for
# every row whose status is "Handback transform failed" or "Ready for handoff"
# (i.e. every row that is still awaiting/being re-handed-off), stamping them
# with the current run time. Rows that are already in sync or in translation
# are left untouched.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$overviewStamp = "2016-26-13 08:26:16"
$zhcnStamp     = "2016-03-13 08:26:12"
$dedeStamp     = "2016-03-13 08:26:16"

# Rows (1-based, header is row 1) whose handoff needs to be (re)stamped.
$rows = @(7, 10, 11, 12, 13, 14, 15, 16)

foreach ($r in $rows) {
    $overview.Range("D$r").Value = $overviewStamp
    $zhcn.Range("E$r").Value = $zhcnStamp
    $dede.Range("E$r").Value = $dedeStamp
}
